# Add 2022-Q4 data
#
# The workbook tracks quarterly "持仓" snapshots, one worksheet per quarter,
# plus a "总计" (totals) summary sheet. This adds a new quarter (2022-Q4):
#   1. A new worksheet named "2022-Q4" is inserted right after "总计" (so all
#      the older quarter sheets shift right by one tab, unchanged otherwise).
#   2. The new sheet is seeded by duplicating the "2022-Q3" sheet (so it
#      picks up the same layout/column styling) and then its data cells are
#      overwritten with the real 2022-Q4 numbers.
#   3. The "总计" summary sheet gets a new row inserted right under its
#      header for the 2022-Q4 totals, and the running index column is
#      renumbered.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. 总计 (totals) sheet: insert a new second row for 2022-Q4.
# ---------------------------------------------------------------------------
$totals = $wb.Worksheets.Item("总计")

$totals.Rows("2:2").Insert()

# New row 2 should look like the other data rows (no bold/border), so copy
# the formatting down from row 3 (which still holds the old row 2 content).
$totals.Range("A3:D3").Copy()
$totals.Range("A2:D2").PasteSpecial(-4122)

$totals.Cells.Item(2, 1).Value = 0
$totals.Cells.Item(2, 2).Value = "2022-Q4"
$totals.Cells.Item(2, 3).Value = 6
$totals.Cells.Item(2, 4).Value = 0.19

# Renumber the 0-based running index in column A for the rest of the rows.
for ($r = 3; $r -le 8; $r++) {
    $totals.Cells.Item($r, 1).Value = $r - 2
}

# ---------------------------------------------------------------------------
# 2. Create the "2022-Q4" worksheet right after "总计" by duplicating the
#    "2022-Q3" sheet (keeps the same header/column styles), then overwrite
#    its contents with the 2022-Q4 fund breakdown.
# ---------------------------------------------------------------------------
$q3 = $wb.Worksheets.Item("2022-Q3")
$q3.Copy($q3)

$q4 = $wb.Worksheets.Item(2)
$q4.Name = "2022-Q4"

$q4Data = @(
    @(0, "000369", "广发全球医疗保健（QDII）人民币A", "3.16", "80.87", "2.77", "0.0875", 7),
    @(1, "000370", "广发全球医疗保健（QDII）美元A",   "3.16", "80.87", "2.77", "0.0875", 7),
    @(2, "014002", "浦银安盛全球智能科技股票（QDII）C", "0.30", "42.55", "1.35", "0.0040", 7),
    @(3, "006555", "浦银安盛全球智能科技股票（QDII）A", "0.25", "42.55", "1.35", "0.0034", 7),
    @(4, "016280", "广发全球医疗保健（QDII）人民币C",   "0.10", "80.87", "2.77", "0.0028", 7),
    @(5, "016281", "广发全球医疗保健（QDII）美元C",     "0.10", "80.87", "2.77", "0.0028", 7)
)

for ($i = 0; $i -lt $q4Data.Length; $i++) {
    $row = $i + 2
    $vals = $q4Data[$i]

    $q4.Cells.Item($row, 1).Value = $vals[0]
    $q4.Cells.Item($row, 2).Value = "'" + $vals[1]
    $q4.Cells.Item($row, 3).Value = "'" + $vals[2]
    $q4.Cells.Item($row, 4).Value = "'" + $vals[3]
    $q4.Cells.Item($row, 5).Value = "'" + $vals[4]
    $q4.Cells.Item($row, 6).Value = "'" + $vals[5]
    $q4.Cells.Item($row, 7).Value = "'" + $vals[6]
    $q4.Cells.Item($row, 8).Value = $vals[7]
}
